$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '20.541.84'
$ws.Range("E2").Value = '  +2.58%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.471.41'
$ws.Range("E3").Value = '  +3.58%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  +0.69%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9440'
$ws.Range("E5").Value = '  -5.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '281.18'
$ws.Range("E6").Value = '  +2.66%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3718'
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3190'
$ws.Range("E8").Value = '  +3.75%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '41.42'
$ws.Range("E9").Value = '  +3.84%  '
$ws.Range("E10").Value = '  +4.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06680'
$ws.Range("E11").Value = '  +1.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.604'
$ws.Range("E13").Value = '  +3.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.32'
$ws.Range("E14").Value = '  +6.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.241'
$ws.Range("E15").Value = '  +1.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.473.82'
$ws.Range("E16").Value = '  +3.80%  '
$ws.Range("E17").Value = '  +2.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9405'
$ws.Range("E18").Value = '  -5.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.05752'
$ws.Range("E19").Value = '  -0.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.32'
$ws.Range("E20").Value = '  -3.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.701'
$ws.Range("E21").Value = '  +1.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.80'
$ws.Range("E22").Value = '  +2.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.26'
$ws.Range("E23").Value = '  +2.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.263'
$ws.Range("E24").Value = '  -2.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '20.833.94'
$ws.Range("E25").Value = '  +4.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.311'
$ws.Range("E26").Value = '  +0.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '138.19'
$ws.Range("E27").Value = '  -0.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.61'
$ws.Range("E28").Value = '  +4.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.638.62'
$ws.Range("E29").Value = '  +3.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '113.75'
$ws.Range("E30").Value = '  +4.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.933'
$ws.Range("E31").Value = '  +3.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.318'
$ws.Range("E32").Value = '  -1.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8526'
$ws.Range("E33").Value = '  -3.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.611'
$ws.Range("E34").Value = '  +26.57%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07843'
$ws.Range("E35").Value = '  +1.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06051'
$ws.Range("E36").Value = '  +5.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.940'
$ws.Range("E37").Value = '  +3.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.73'
$ws.Range("E38").Value = '  -5.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02076'
$ws.Range("E39").Value = '  +1.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.121'
$ws.Range("E40").Value = '  +2.99%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.631'
$ws.Range("E41").Value = '  -9.68%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1907'
$ws.Range("E42").Value = '  -1.04%  '
$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9575'
$ws.Range("E43").Value = '  -4.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5420'
$ws.Range("E44").Value = '  +1.84%  '
$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.591'
$ws.Range("E45").Value = '  +1.52%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.47'
$ws.Range("E46").Value = '  +1.85%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.35'
$ws.Range("E47").Value = '  +11.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5337'
$ws.Range("E48").Value = '  +4.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.831'
$ws.Range("E49").Value = '  +1.83%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06456'
$ws.Range("E50").Value = '  +4.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.048'